$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 271
$ws.Range("F4").Value = 1104
$ws.Range("F5").Value = 2645
$ws.Range("F7").Value = 683
$ws.Range("F8").Value = 59
$ws.Range("F9").Value = 249
$ws.Range("F11").Value = 695
$ws.Range("F12").Value = 93
$ws.Range("F13").Value = 121
$ws.Range("F14").Value = 1542
$ws.Range("F15").Value = 303
$ws.Range("F17").Value = 196

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 25
$ws.Range("F7").Value = 21
$ws.Range("F12").Value = 44

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6349
$ws.Range("F3").Value = 794
$ws.Range("F4").Value = 2018

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6349
$ws.Range("F3").Value = 794
$ws.Range("F4").Value = 2018
$ws.Range("F9").Value = 25
$ws.Range("F11").Value = 271
$ws.Range("F12").Value = 1104
$ws.Range("F14").Value = 21
$ws.Range("F16").Value = 2645
$ws.Range("F21").Value = 44
$ws.Range("F22").Value = 683
$ws.Range("F23").Value = 59
$ws.Range("F24").Value = 249
$ws.Range("F27").Value = 695
$ws.Range("F28").Value = 93
$ws.Range("F29").Value = 121
$ws.Range("F31").Value = 1542
$ws.Range("F32").Value = 303
$ws.Range("F36").Value = 196

